# Corrigindo exibicao exame imagem
$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("OTAVIO RAMOS DE ALMEIDA   ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Text = "FABRICIO SANCHEZ BERGAMIN   "
} else {
    Write-Output "NOT FOUND: OTAVIO RAMOS DE ALMEIDA   "
}

$r = $d.Content
$r.Find.Execute("15/02/1988   ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Text = "10/10/1981   "
} else {
    Write-Output "NOT FOUND: 15/02/1988   "
}

$r = $d.Content
$r.Find.Execute("26294", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Text = "576"
} else {
    Write-Output "NOT FOUND: 26294"
}

$r = $d.Content
$r.Find.Execute("MARIA IRACY RAMOS DOS SANTOS   ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Text = "JUSSARA SANCHEZ BERGAMIN   "
} else {
    Write-Output "NOT FOUND: MARIA IRACY RAMOS DOS SANTOS   "
}

$r = $d.Content
$r.Find.Execute("20/02/2019   ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Text = "22/10/2019   "
} else {
    Write-Output "NOT FOUND: 20/02/2019   "
}

$oldReport = @'
.
TÉCNICA
Exame  realizado  com cortes tomográficos computadorizados axiais,  sem a infusão endovenosa de contraste iodado, segundo solicitação do médico assistente.
Salientamos que a não utilização do meio de contraste iodado por via endovenosa prejudica a adequada caracterização das estruturas abdominais.
ANÁLISE
Fígado de topografia, morfologia, situação e dimensões, preservadas, exibindo coeficientes de atenuação homogêneos.
Não há evidência de dilatação das vias biliares intra ou extra-hepáticas, bem como da vesícula biliar.
Baço, pâncreas e adrenais com topografia, dimensões, contornos e densidade normais.
Rins de topografia, morfologia e dimensões preservadas, com coeficientes de atenuação homogêneos, sem a caracterização de hidronefrose.
Dois cálculos não obstrutivos no terço médio e inferior do rim esquerdo medindo até 0,4 cm.
Aorta e veia cava inferior com posição e calibre normais.
Ausência de linfonodomegalias, líquido livre ou de coleções organizadas na cavidade abdominal.
Bexiga urinária em pequena repleção, com paredes lisas e regulares e conteúdo homogêneo.
Próstata e vesículas seminais sem alterações detectáveis ao método.
OPINIÃO
Nefrolitíase esquerda não obstrutiva.
           Dra. Amanda Prist
             CRM-MG: 56.487

'@
$newReport = @'
{\rtf1\ansi\ansicpg1252\deff0\deflang1046{\fonttbl{\f0\fnil\fcharset0 }{\f1\fswiss\fcharset0  }{\f2\fswiss\fprq2\fcharset0  }}
{\colortbl ;\red0\green0\blue0;}
\viewkind4\uc1\pard\f0\fs16 .\par
\par
\par
ECOCARDIOGRAMA\par
\par
\par
\cf1\b\f1\fs20  \par
\pard\qj\cf0\b0\f2 R\'cdTMO:\par
Paciente em ritmo card\'edaco regular.\par
\par
C\'c2MARAS CARD\'cdACAS:\par
C\'e2maras card\'edacas com dimens\'f5es normais.\par
\par
VENTR\'cdCULOS:\par
Ventr\'edculo esquerdo apresenta espessura e fun\'e7\'e3o sist\'f3lica preservadas, n\'e3o sendo observadas altera\'e7\'f5es da contra\'e7\'e3o segmentar de parede.\par
An\'e1lise da fun\'e7\'e3o diast\'f3lica do ventr\'edculo esquerdo com padr\'e3o normal.\par
Ventr\'edculo direito apresenta fun\'e7\'e3o sistolica dento da normalidade,\par
\par
V\'c1LVULA MITRAL:\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas cuspides.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
V\'c1LVULA A\'d3RTICA:\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas valvulas.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
V\'c1LVULA TRIC\'daSPIDE;\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas c\'faspides.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
VALVA PULMONAR:\par
Apresenta aspecto e movimenta\'e7\'e3o normais de suas v\'e1lvulas.\par
O estudo com Doppler e mapeamento de fluxo em cores s\'e3o normais.\par
\par
PERIC\'c1RDIO:\par
Peric\'e1rdio com aspecto ecocardiografico normal.\par
\par
AORTA:\par
Seios a\'f3rticos, aorta ascendente e arco a\'f3rtico com dimens\'f5es e fluxos normais.\par
\par
IMPRESS\'c3O DIAGN\'d3STICA:\par
_________________________ \par
\par
\pard ECODOPPLERCARDIOGRAMA DENTRO DA NORMALIDADE PARA O BIOTIPO E FAIXA ET\'c1RIA.\f0\fs16\par
}
_x0000_
'@

$r = $d.Content
$r.Find.Execute($oldReport, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Text = $newReport
    Write-Output "Report text replaced successfully"
} else {
    Write-Output "NOT FOUND: report text block"
}
